# Update the "想去人数" (want-to-go count) column F values on the
# "展览" and "全部类型" worksheets, per the output regenerated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> hashtable of (row -> new F value)
$updates = @{
    "展览" = @{
        2  = 210
        3  = 5574
        4  = 39
        5  = 63
        7  = 667
        8  = 650
        10 = 3
        11 = 1083
        13 = 1561
        14 = 5150
        16 = 261
        17 = 227
        18 = 38
        20 = 115
        21 = 4448
        22 = 220
        23 = 1173
        24 = 122
        25 = 70
        26 = 213
        27 = 60
        28 = 182
        30 = 151
        32 = 349
        35 = 68
        36 = 11
        38 = 2
    }
    "全部类型" = @{
        2  = 210
        4  = 5574
        5  = 39
        6  = 63
        8  = 667
        9  = 650
        11 = 3
        12 = 1083
        14 = 1561
        15 = 5150
        17 = 261
        18 = 227
        19 = 38
        21 = 115
        22 = 4448
        23 = 220
        24 = 1173
        25 = 122
        26 = 70
        27 = 213
        28 = 60
        29 = 182
        31 = 151
        33 = 349
        36 = 68
        37 = 11
        39 = 2
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
